# Update the Pomc-Mc4r LR-pair sheet with newly computed TPM-based values.
# Row 2 changes from ECs -> ECs to ECs -> MuSCs (target cluster updated).
# Row 3 changes from ECs -> MuSCs to MuSCs -> MuSCs (sending cluster updated).
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the various expression / specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row2 = @(
    "ECs", "Pomc", "Mc4r", "MuSCs",
    3, 1, 0.5131323333333334, 1.539397,
    0.9119558630037493, 0.9119558630037494,
    3, 1, 0.309635, 0.928905,
    1, 1, 0.1588837300316667, 1.429953570285,
    0.9119558630037493, 0.9119558630037494
)

$row3 = @(
    "MuSCs", "Pomc", "Mc4r", "MuSCs",
    1, 0.3333333333333333, 0.04954, 0.14862,
    0.08804413699625062, 0.08804413699625063,
    3, 1, 0.309635, 0.928905,
    1, 1, 0.0153393179, 0.1380538611,
    0.08804413699625062, 0.08804413699625063
)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
